$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.204.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.06%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.791.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.48%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '344.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.31%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.43%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.549'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.38%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.579'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.11%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.37%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0855'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.32%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.15'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.64%  '

# Row 13
$ws.Range('E13').Value = '  +2.23%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.15%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.220.35'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.29%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.776.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.80%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.885'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.56%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '52.004.72'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.06%  '

# Row 19
$ws.Range('E19').Value = '  +8.18%  '

# Row 20
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.87%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.23%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0980'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.05%  '

# Row 23
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.24%  '

# Row 24
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.28%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.61%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.97%  '

# Row 27
$ws.Range('E27').Value = '  -0.11%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.77%  '

# Row 29
$ws.Range('E29').Value = '  +0.39%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.140'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.58%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.68%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.44%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.26%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0823'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.32%  '

# Row 35
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0410'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +17.20%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.11%  '

# Row 37
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.10'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.61%  '

# Row 38
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.02'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.36%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.96'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.32%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.23'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.41%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +24.06%  '

# Row 42
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.116'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.74%  '

# Row 43
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.07%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.99'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.39%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.32'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.53%  '

# Row 46
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.34'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.15%  '

# Row 47
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.073.83'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.79%  '

# Row 48
$ws.Range('E48').Value = '  +1.60%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.25%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.905'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.02%  '
